$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (prices with "." as thousands separators,
# leading zeros, trailing zeros, etc.) that must stay literal text, so force the
# Text number format before writing those values (mirrors typing into a text cell).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.041.78"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.845.07"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.36"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  +2.94%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.51"
$ws.Range("E8").Value = "  +5.52%  "
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("E11").Value = "  -1.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.114.99"
$ws.Range("E12").Value = "  +2.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.52"
$ws.Range("E13").Value = "  +4.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.850.36"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.674"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.69"
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.088.13"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.08"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0790"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.70"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.18"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.77"
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +3.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.86"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.89"
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("E27").Value = "  +2.08%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.124"
$ws.Range("E28").Value = "  +3.96%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.74"
$ws.Range("E29").Value = "  +12.42%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.65"
$ws.Range("E34").Value = "  +23.60%  "
$ws.Range("E35").Value = "  +10.87%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.760"
$ws.Range("E36").Value = "  +8.77%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.25"
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("E38").Value = "  +10.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0201"
$ws.Range("E39").Value = "  +4.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "90.05"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.345.77"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.54"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("E43").Value = "  +3.09%  "
$ws.Range("E44").Value = "  +4.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.87"
$ws.Range("E46").Value = "  +78.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0531"
$ws.Range("E47").Value = "  +3.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.33"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.033.88"
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E50").Value = "  +15.33%  "
$ws.Range("E51").Value = "  +0.79%  "
